$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update scattered F-column values (rows 5, 8, 12, 14, 18) ---
$ws.Range("F5").Value = ""
$ws.Range("F8").Value = 17.05
$ws.Range("F12").Value = ""
$ws.Range("F14").Value = 17.76
$ws.Range("F18").Value = ""

# --- Remove the "SC 92" row (originally row 28) and the "RM 232" row (originally row 26) ---
# Delete the lower row index first so the upper one keeps its original position.
$ws.Rows(28).Delete()
$ws.Rows(26).Delete()

# --- Fill in / clear a few individual cells on the now-shifted rows ---
# Row 26 now holds what was "SC 5" -> column C (B) becomes 10.8
$ws.Range("C26").Value = 10.8
# Row 27 now holds what was "SC 101" -> column C (B) becomes blank
$ws.Range("C27").Value = ""
# Row 33 now holds what was "SC 232" -> column D (C) becomes -14.1
$ws.Range("D33").Value = -14.1
